$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking values so Excel does not
# silently convert them to numbers (matches original inline-string cells).
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '43.672.22'
Set-TextValue $ws.Range("E2") '  -0.52%  '

# Row 3
Set-TextValue $ws.Range("D3") '2.292.64'
Set-TextValue $ws.Range("E3") '  +0.08%  '

# Row 4
Set-TextValue $ws.Range("E4") '  -0.26%  '

# Row 5
Set-TextValue $ws.Range("D5") '117.59'
Set-TextValue $ws.Range("E5") '  +3.72%  '

# Row 6
Set-TextValue $ws.Range("D6") '267.89'
Set-TextValue $ws.Range("E6") '  -0.96%  '

# Row 7
Set-TextValue $ws.Range("E7") '  +2.99%  '

# Row 8
Set-TextValue $ws.Range("E8") '  -0.03%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.618'
Set-TextValue $ws.Range("E9") '  -0.32%  '

# Row 10
Set-TextValue $ws.Range("D10") '48.31'
Set-TextValue $ws.Range("E10") '  +1.19%  '

# Row 11
Set-TextValue $ws.Range("D11") '0.0942'
Set-TextValue $ws.Range("E11") '  -0.33%  '

# Row 12
Set-TextValue $ws.Range("D12") '9.23'
Set-TextValue $ws.Range("E12") '  +1.89%  '

# Row 13
Set-TextValue $ws.Range("E13") '  +1.14%  '

# Row 14
Set-TextValue $ws.Range("D14") '15.57'
Set-TextValue $ws.Range("E14") '  -1.62%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.888'
Set-TextValue $ws.Range("E15") '  +4.42%  '

# Row 16
Set-TextValue $ws.Range("D16") '2.633.85'
Set-TextValue $ws.Range("E16") '  -0.16%  '

# Row 17
Set-TextValue $ws.Range("D17") '2.287.77'
Set-TextValue $ws.Range("E17") '  -0.96%  '

# Row 18
Set-TextValue $ws.Range("D18") '43.738.58'
Set-TextValue $ws.Range("E18") '  -0.03%  '

# Row 19
Set-TextValue $ws.Range("E19") '  +0.03%  '

# Row 20
Set-TextValue $ws.Range("D20") '6.99'
Set-TextValue $ws.Range("E20") '  +2.68%  '

# Row 21
Set-TextValue $ws.Range("D21") '72.59'
Set-TextValue $ws.Range("E21") '  +0.20%  '

# Row 22
Set-TextValue $ws.Range("D22") '2.48'
Set-TextValue $ws.Range("E22") '  +1.91%  '

# Row 23
Set-TextValue $ws.Range("D23") '236.64'
Set-TextValue $ws.Range("E23") '  +1.73%  '

# Row 24
Set-TextValue $ws.Range("D24") '9.74'
Set-TextValue $ws.Range("E24") '  +0.73%  '

# Row 25
Set-TextValue $ws.Range("D25") '2.90'
Set-TextValue $ws.Range("E25") '  -3.05%  '

# Row 26
Set-TextValue $ws.Range("E26") '  +1.83%  '

# Row 27
Set-TextValue $ws.Range("D27") '11.87'
Set-TextValue $ws.Range("E27") '  +1.88%  '

# Row 28
Set-TextValue $ws.Range("D28") '42.66'
Set-TextValue $ws.Range("E28") '  +2.30%  '

# Row 29
Set-TextValue $ws.Range("D29") '3.41'
Set-TextValue $ws.Range("E29") '  +0.50%  '

# Row 30
Set-TextValue $ws.Range("E30") '  -0.29%  '

# Row 31
Set-TextValue $ws.Range("D31") '174.08'
Set-TextValue $ws.Range("E31") '  -0.58%  '

# Row 32
Set-TextValue $ws.Range("D32") '21.86'
Set-TextValue $ws.Range("E32") '  +1.36%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.0917'
Set-TextValue $ws.Range("E33") '  -1.53%  '

# Row 34
Set-TextValue $ws.Range("D34") '5.76'
Set-TextValue $ws.Range("E34") '  +1.37%  '

# Row 35
Set-TextValue $ws.Range("E35") '  +2.62%  '

# Row 36
Set-TextValue $ws.Range("D36") '0.0384'
Set-TextValue $ws.Range("E36") '  +5.75%  '

# Row 37
Set-TextValue $ws.Range("D37") '4.72'
Set-TextValue $ws.Range("E37") '  +1.56%  '

# Row 38
Set-TextValue $ws.Range("D38") '3.97'
Set-TextValue $ws.Range("E38") '  +4.96%  '

# Row 39
Set-TextValue $ws.Range("E39") '  +0.17%  '

# Row 40
Set-TextValue $ws.Range("D40") '2.58'
Set-TextValue $ws.Range("E40") '  +8.83%  '

# Row 41
Set-TextValue $ws.Range("D41") '14.28'
Set-TextValue $ws.Range("E41") '  +5.13%  '

# Row 42
Set-TextValue $ws.Range("D42") '74.65'
Set-TextValue $ws.Range("E42") '  +1.02%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.238'
Set-TextValue $ws.Range("E43") '  -1.97%  '

# Row 44
Set-TextValue $ws.Range("E44") '  -5.65%  '

# Row 45
Set-TextValue $ws.Range("E45") '  -0.23%  '

# Row 46
Set-TextValue $ws.Range("E46") '  -1.24%  '

# Row 47
Set-TextValue $ws.Range("D47") '1.29'
Set-TextValue $ws.Range("E47") '  +4.37%  '

# Row 48
Set-TextValue $ws.Range("E48") '  -2.03%  '

# Row 49
Set-TextValue $ws.Range("B49") 'ordi'
Set-TextValue $ws.Range("C49") 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
Set-TextValue $ws.Range("D49") '73.58'
Set-TextValue $ws.Range("E49") '  +37.34%  '

# Row 50
Set-TextValue $ws.Range("B50") 'Cronos'
Set-TextValue $ws.Range("C50") 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.101'
Set-TextValue $ws.Range("E50") '  +0.75%  '

# Row 51
Set-TextValue $ws.Range("D51") '101.99'
Set-TextValue $ws.Range("E51") '  +0.34%  '
